# C5-PowerPoint.pptx edit
# Commit: Tue, Jul 07, 2020 10:06:08 AM
#
# The recorded change swaps the table style applied to the "Sources of
# Finance" table (slide 6, shape 2) away from the deck's custom
# "Table_0" style and onto a different (built-in) table style id.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Table styles can't be set through the .Style property directly -
# PowerPoint exposes this as a method call.
$tbl.ApplyStyle("{8B266275-4925-44DE-80B1-7EBA46EFA41B}")
